$d = $word.ActiveDocument

# --- Paragraph 1 ("Sentri:") ---
# Replace centered alignment with a pair of tab stops (center @ 4513, left @ 7536)
$p1 = $d.Paragraphs.Item(1)
$p1.Range.ParagraphFormat.Alignment = 0
$p1.Range.ParagraphFormat.TabStops.Add(4513 / 20, 1)
$p1.Range.ParagraphFormat.TabStops.Add(7536 / 20, 0)

# Insert a tab before "Sentri" and a tab after ":" (text-wise); formatting is
# fixed up below via the WordOpenXML round trip so each tab becomes its own
# <w:tab/> run instead of a literal tab character merged into a text run.
$rng1 = $d.Range($p1.Range.Start, $p1.Range.End)
$rng1.Find.Execute("Sentri", $false, $false, $false, $false, $false, $true, 1, $false, "`tSentri", 2)

$p1b = $d.Paragraphs.Item(1)
$rng2 = $p1b.Range
$rng2.Find.Execute(":", $false, $false, $false, $false, $false, $true, 1, $false, ":`t", 2)

# --- Paragraph 3 ("Authors:") ---
$p3 = $d.Paragraphs.Item(3)
$p3.Range.ParagraphFormat.TabStops.Add(6696 / 20, 0)
$p3r = $p3.Range
$p3r.Find.Execute("Authors:", $false, $false, $false, $false, $false, $true, 1, $false, "Authors:`t", 2)

# --- Fix up the inserted tab characters so they serialize as real <w:tab/>
#     elements in their own run (matching how Word represents Tab keypresses)
#     instead of a literal U+0009 inside a <w:t> text run. ---
$xml = $d.Content.WordOpenXML

$xml = $xml -replace '<w:t xml:space="preserve">\tSentri</w:t>', '<w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>Sentri</w:t>'
$xml = $xml -replace '<w:t xml:space="preserve">:\t</w:t>', '<w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:tab/>'
$xml = $xml -replace '<w:t xml:space="preserve">Authors:\t</w:t>', '<w:t>Authors:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:tab/>'

$d.Content.WordOpenXML = $xml
